# Convert both "Platform Coverage" and "MarketShare" sheets from yearly
# (2018-2040) to half-yearly (2018, 2018.5, 2019, 2018.5, ..., 2040) columns,
# and rewrite the coverage / market-share figures so that the "scenario 2"
# trichuris coverage switches to a 6-month (half-yearly) reporting cadence
# starting in 2026.

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet 1: "Platform Coverage"
# -------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# --- Row 1: header years, half-yearly from col H (8) to col AZ (52) ---
$year = 2018.0
for ($col = 8; $col -le 52; $col++) {
    $ws1.Cells.Item(1, $col).Value = $year
    $year = $year + 0.5
}

# --- Row 2 ("All"/Treatment/Campaign/MDA, age 5-15): 0.6 coverage,
#     now also populated for whole years 2022-2025 (cols P,R,T,V) ---
$ws1.Cells.Item(2, 16).Value = 0.6   # P2 = 2022
$ws1.Cells.Item(2, 18).Value = 0.6   # R2 = 2023
$ws1.Cells.Item(2, 20).Value = 0.6   # T2 = 2024
$ws1.Cells.Item(2, 22).Value = 0.6   # V2 = 2025

# --- Rows 3,4,5: coverage values that used to be entered every other
#     year from col P onward now move two columns to the right (since
#     every year now spans two columns) AND become continuous half-year
#     entries from 2026.0 onward (col X = 24) through 2040.0 (col AZ = 52).
#     First clear the stale entries at the old (now mis-aligned) columns. ---
$oldCols = 16,18,20,22   # P,R,T,V
foreach ($r in 3,4,5) {
    foreach ($c in $oldCols) {
        $ws1.Cells.Item($r, $c).ClearContents()
    }
}

$rowValues = @{ 3 = 0.8; 4 = 0.5; 5 = 0.5 }
foreach ($r in 3,4,5) {
    $val = $rowValues[$r]
    for ($col = 24; $col -le 52; $col++) {
        $ws1.Cells.Item($r, $col).Value = $val
    }
}

# --- View state: selection / scroll position ---
$ws1.Activate()
[void]$ws1.Range("AZ17").Select()

# -------------------------------------------------------------------------
# Sheet 2: "MarketShare"
# -------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# --- Row 1: header years, half-yearly from col D (4) to col AV (48) ---
$year = 2018.0
for ($col = 4; $col -le 48; $col++) {
    $ws2.Cells.Item(1, $col).Value = $year
    $year = $year + 0.5
}

# --- Row 2 (DRC/MDA/New Product A): market share of 1, used to start at
#     whole-year col L (2026) through Z (2040); now starts at half-year
#     col T (2026.0) through AV (2040.0). Clear the stale pre-2026 whole
#     year entries that no longer land on 2026+ under the new mapping. ---
foreach ($c in 12,13,14,15,16,17,18,19) {   # L..S (old 2026-2033 whole years)
    $ws2.Cells.Item(2, $c).ClearContents()
}
for ($col = 20; $col -le 48; $col++) {       # T..AV (2026.0 .. 2040.0)
    $ws2.Cells.Item(2, $col).Value = 1
}

# --- Row 3 (DRC/MDA/Old Product B (SOC)): market share of 1, used to run
#     whole years 2018-2025 (cols D..K); now fills every half-year through
#     2025.5 (cols D..S), since New Product A takes over entirely at 2026. ---
for ($col = 4; $col -le 19; $col++) {        # D..S (2018.0 .. 2025.5)
    $ws2.Cells.Item(3, $col).Value = 1
}

# --- View state: selection / scroll position ---
$ws2.Activate()
[void]$ws2.Range("Q3").Select()

# Leave MarketShare as the active sheet (matches the original workbook).
$ws2.Activate()
